$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @('1.02','1.081497900768994','1.086543546774944','1.094684272638825','1.100593467886747','1.069232182326978','1.086372186377141','1.089201118891153','1.097320991355226','1.10321525773721','1.032057703193815')
    3  = @('1.02','1.082723370382783','1.087562913519939','1.095860104372401','1.101756183151815','1.069697596471883','1.087258106054709','1.090039401052493','1.098316874825576','1.104199136042878','1.032379641402154')
    4  = @('1.02','1.083516028141429','1.08822223465871','1.096620950462821','1.102508551503431','1.069997363070816','1.087830471195232','1.090580933083967','1.098960695597579','1.104835194292435','1.032587109034625')
    5  = @('1.02','1.0838491909334','1.088499347971876','1.0969408126194','1.102824851800687','1.07012305355252','1.08807088311872','1.090808380592364','1.099231219744461','1.105102456354237','1.032674125760895')
    6  = @('1.02','1.083905126354892','1.088545872719489','1.096994519068599','1.102877960323347','1.070144138126552','1.088111237038926','1.090846557598167','1.09927663387914','1.105147322843559','1.032688724377186')
    7  = @('1.02','1.08352048015356','1.088225937715548','1.096625224465155','1.102512777903567','1.069999043854382','1.087833684416166','1.090583973082497','1.098964310895461','1.104838765997235','1.032588272552303')
    8  = @('1.02','1.081912117918906','1.086888104616937','1.095081650437757','1.100986410893614','1.069389758963919','1.086671770916189','1.089484605766429','1.097657676231529','1.103547884207909','1.032166679481454')
    9  = @('1.02','1.079075566976078','1.084528498518432','1.092361644951289','1.09829680310272','1.068305460229577','1.084617501903155','1.087540493203887','1.095350698575415','1.10126871113912','1.031417269860147')
    10 = @('1.02','1.077182783388337','1.082953895077291','1.090548195233752','1.096503675668195','1.067575377043979','1.083243315288564','1.086239714764203','1.093809588438083','1.099746172413674','1.030913265678515')
    11 = @('1.02','1.076362742944861','1.082271694166328','1.089762903493247','1.095727198474868','1.06725751964304','1.082647149710902','1.085675330850244','1.093141509732641','1.099086144537673','1.030693978060254')
    12 = @('1.02','1.076058072788601','1.082018234659152','1.089471200825488','1.095438772259455','1.067139192798069','1.082425535215434','1.085465520976034','1.092893238394758','1.098840865107971','1.030612366557666')
    13 = @('1.02','1.076123428809248','1.082072605310734','1.089533772537336','1.095500641053731','1.067164586111043','1.082473080090677','1.085510533739829','1.092946498770213','1.098893483641141','1.030629879674144')
    14 = @('1.02','1.076337560270535','1.082250744351443','1.089738791482048','1.095703357247967','1.067247744028411','1.082628834507464','1.085657991429709','1.093120989946861','1.099065872033493','1.030687235260476')
    15 = @('1.02','1.07646948429449','1.082360493735053','1.089865108991069','1.095828256269204','1.067298945821599','1.082724777093387','1.085748822067032','1.093228484140934','1.099172070830058','1.030722552936813')
    16 = @('1.02','1.077237197011577','1.082999162215392','1.090600311114806','1.09655520687883','1.067596435725921','1.083282856801114','1.086277146951813','1.093853910290267','1.099789960161002','1.030927796900892')
    17 = @('1.02','1.077718639879817','1.083399677322884','1.09106146761622','1.097011191476476','1.067782580286827','1.083632620711793','1.086608245355415','1.094246016745376','1.100177341509216','1.03105625942404')
    18 = @('1.02','1.077999413862455','1.083633253751978','1.091330447125512','1.097277155729022','1.067890988743562','1.083836522815143','1.086801259732008','1.094474651843537','1.100403221423226','1.031131088117986')
    19 = @('1.02','1.078095143235458','1.083712891004701','1.091422161350094','1.097367842067909','1.067927925013442','1.083906029681178','1.086867054082386','1.094552598012092','1.100480228293953','1.031156585597611')
    20 = @('1.02','1.077666990135624','1.083356709708025','1.091011990448717','1.096962269056906','1.067762625968887','1.083595105669867','1.086572732991167','1.094203955097975','1.100135786728361','1.031042487095256')
    21 = @('1.02','1.076274505862849','1.082198288504712','1.08967841881403','1.095643662617223','1.067223263292691','1.082582973477276','1.085614573591254','1.093069609893093','1.09901511115605','1.030670349844157')
    22 = @('1.02','1.075398585830993','1.08146959686682','1.08883988855651','1.094814555144495','1.066882637142659','1.081945609477615','1.085011141603768','1.092355723292922','1.098309827748804','1.030435456163999')
    23 = @('1.02','1.075862967341553','1.081855923426293','1.089284415724281','1.095254085794607','1.067063352749648','1.082283583210899','1.085331127542619','1.092734233026425','1.098683775930756','1.030560064760422')
    24 = @('1.02','1.077690328566203','1.083376125038869','1.091034347067374','1.096984375004449','1.067771642980774','1.083612057440163','1.086588779839031','1.094222961174507','1.100154563770972','1.031048710530929')
    25 = @('1.02','1.079809182347173','1.085138778148838','1.093064844683461','1.098992134634917','1.068587045918227','1.085149397442902','1.088043916546356','1.095947652113215','1.101858470650579','1.031611783761389')
}

$cols = @(2,3,4,5,6,9,10,11,12,13,14)  # B,C,D,E,F,I,J,K,L,M,N

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = [double]$vals[$i]
    }
}
